# Error Calculations and Plots
# Apply edits to the "missing data" worksheet:
#  - A few cells switch between having a numeric value and being blank
#    (simulating different randomly removed / restored data points).
#  - The row for "RM 232" is removed entirely.
#  - The row for "SC 92" is removed entirely.
#  - The (now shifted) row for "SC 193" gets a value restored in column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Individual cell value changes (using original row numbering) ---

# Row 3 (RM 8): E3 goes from blank to -5.7
$ws.Range("E3").Value = -5.7

# Row 4 (RM 9): F4 goes from 17.97 to blank
$ws.Range("F4").ClearContents()

# Row 5 (RM 14): E5 goes from -5 to blank
$ws.Range("E5").ClearContents()

# Row 9 (RM 42): F9 goes from blank to 17.26
$ws.Range("F9").Value = 17.26

# Row 10 (RM 52 a): F10 goes from blank to 16.43
$ws.Range("F10").Value = 16.43

# Row 13 (RM 88): F13 goes from 17.1 to blank
$ws.Range("F13").ClearContents()

# Row 14 (RM 90): F14 goes from 17.76 to blank
$ws.Range("F14").ClearContents()

# Row 21 (RM 135): E21 goes from blank to -8.699999999999999
$ws.Range("E21").Value = -8.699999999999999

# Row 23 (RM 140): E23 goes from -7 to blank
$ws.Range("E23").ClearContents()

# --- Remove the "RM 232" row entirely ---
$ws.Rows(26).Delete()

# After removing RM 232, "SC 92" (originally row 28) is now row 27
$ws.Rows(27).Delete()

# After both deletions, "SC 193" (originally row 34) is now row 32.
# Restore its column E value.
$ws.Range("E32").Value = -6.4
